$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 9, shifting existing rows 9.. down by one.
$ws.Rows.Item(9).EntireRow.Insert()

# Populate the new row 9 with the new market observation.
$ws.Cells.Item(9, 1).Value = 11
$ws.Cells.Item(9, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(9, 3).Value = "Bíobío"
$ws.Cells.Item(9, 4).Value = 44530
$ws.Cells.Item(9, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(9, 5).Value = 8
$ws.Cells.Item(9, 6).Value = 100112043
$ws.Cells.Item(9, 7).Value = "Pepino ensalada"
$ws.Cells.Item(9, 8).Value = "Sin especificar"
$ws.Cells.Item(9, 9).Value = "Primera"
$ws.Cells.Item(9, 10).Value = 350
$ws.Cells.Item(9, 11).Value = 6000
$ws.Cells.Item(9, 12).Value = 6500
$ws.Cells.Item(9, 13).Value = 6286
$ws.Cells.Item(9, 14).Value = "`$/caja 80 unidades"
$ws.Cells.Item(9, 15).Value = "Región del Maule"
$ws.Cells.Item(9, 16).Value = 79
$ws.Cells.Item(9, 17).Value = 80
$ws.Cells.Item(9, 18).Value = "Hortaliza"
